$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Phase 1: set idx, type, title, image for each new row (row order)
$ws.Range("A29").Value = 28
$ws.Range("B29").Value = 'planet'
$ws.Range("D29").Value = 'Planet Acere'
$ws.Range("I29").Value = 'Planet Acere.png'

$ws.Range("A30").Value = 29
$ws.Range("B30").Value = 'planet'
$ws.Range("D30").Value = 'Planet Carbonic'
$ws.Range("I30").Value = 'Planet Carbonic.png'

$ws.Range("A31").Value = 30
$ws.Range("B31").Value = 'planet'
$ws.Range("D31").Value = 'Planet East Eggor'
$ws.Range("I31").Value = 'Planet East Eggor.png'

$ws.Range("A32").Value = 31
$ws.Range("B32").Value = 'planet'
$ws.Range("D32").Value = 'Planet Eden'
$ws.Range("I32").Value = 'Planet Eden.png'

$ws.Range("A33").Value = 32
$ws.Range("B33").Value = 'planet'
$ws.Range("D33").Value = 'Planet Fetlink'
$ws.Range("I33").Value = 'Planet Fetlink.png'

$ws.Range("A34").Value = 33
$ws.Range("B34").Value = 'planet'
$ws.Range("D34").Value = 'Planet Garbonic'
$ws.Range("I34").Value = 'Planet Garbonic.png'

$ws.Range("A35").Value = 34
$ws.Range("B35").Value = 'planet'
$ws.Range("D35").Value = 'Planet Hearon'
$ws.Range("I35").Value = 'Planet Hearon.png'

$ws.Range("A36").Value = 35
$ws.Range("B36").Value = 'planet'
$ws.Range("D36").Value = 'Planet Hogsworth'
$ws.Range("I36").Value = 'Planet Hogsworth.png'

$ws.Range("A37").Value = 36
$ws.Range("B37").Value = 'planet'
$ws.Range("D37").Value = 'Planet Ikeness'
$ws.Range("I37").Value = 'Planet Ikeness.png'

$ws.Range("A38").Value = 37
$ws.Range("B38").Value = 'planet'
$ws.Range("D38").Value = 'Planet Mitter'
$ws.Range("I38").Value = 'Planet Mitter.png'

$ws.Range("A39").Value = 38
$ws.Range("B39").Value = 'planet'
$ws.Range("D39").Value = 'Planet Napaul'
$ws.Range("I39").Value = 'Planet Napaul.png'

$ws.Range("A40").Value = 39
$ws.Range("B40").Value = 'planet'
$ws.Range("D40").Value = 'Planet Networth'
$ws.Range("I40").Value = 'Planet Networth.png'

$ws.Range("A41").Value = 40
$ws.Range("B41").Value = 'planet'
$ws.Range("D41").Value = 'Planet Pluton'
$ws.Range("I41").Value = 'Planet Pluton.png'

$ws.Range("A42").Value = 41
$ws.Range("B42").Value = 'planet'
$ws.Range("D42").Value = 'Planet Satiron'
$ws.Range("I42").Value = 'Planet Satiron.png'

$ws.Range("A43").Value = 42
$ws.Range("B43").Value = 'planet'
$ws.Range("D43").Value = 'Planet Scysm'
$ws.Range("I43").Value = 'Planet Scysm.png'

$ws.Range("A44").Value = 43
$ws.Range("B44").Value = 'planet'
$ws.Range("D44").Value = 'Planet Starstir'
$ws.Range("I44").Value = 'Planet Starstir.png'

$ws.Range("A45").Value = 44
$ws.Range("B45").Value = 'planet'
$ws.Range("D45").Value = 'Planet Telepan'
$ws.Range("I45").Value = 'Planet Telepan.png'

$ws.Range("A46").Value = 45
$ws.Range("B46").Value = 'planet'
$ws.Range("D46").Value = 'Planet Zastron'
$ws.Range("I46").Value = 'Planet Zastron.png'

$ws.Range("A47").Value = 46
$ws.Range("B47").Value = 'planet'
$ws.Range("D47").Value = 'Planet Zerox'
$ws.Range("I47").Value = 'Planet Zerox.png'

# Phase 2: set descriptions in the order: row35 first, then 29,30,31,(subtype C31),32..47
$ws.Range("E35").Value = 'Hearon this planet is life.'
$ws.Range("E29").Value = 'I looks pink from space, but on the surface it is all red'
$ws.Range("E30").Value = 'Look at all the carbon in the atmosphere'
$ws.Range("E31").Value = 'I hope there is water somewhere on the surface'
$ws.Range("C31").Value = 'Secret weapon'
$ws.Range("E32").Value = 'The perfect place to add to the empire''s trophies'
$ws.Range("E33").Value = 'Mostly covered in lava. Just be careful'
$ws.Range("E34").Value = 'There is snow all over this planet'
$ws.Range("E36").Value = 'A very old planet, but it will do the job'
$ws.Range("E37").Value = 'Named for its ink like look'
$ws.Range("E38").Value = 'Comes with a built in belt'
$ws.Range("E39").Value = 'Has excess of oxygen everywhere'
$ws.Range("E40").Value = 'Has a built in planetary defence system'
$ws.Range("E41").Value = 'Rumored to have Plutonium at its core'
$ws.Range("E42").Value = 'Semi-transparent planet'
$ws.Range("E43").Value = 'It has its own moon, have fun'
$ws.Range("E44").Value = 'Glows from space. Hopefully it is  gold'
$ws.Range("E45").Value = 'A few gas storms, but nothing our shields can''t handle'
$ws.Range("E46").Value = 'Just a big rock. Hopefully the core is usable'
$ws.Range("E47").Value = 'A very nice planet'

# Update the active selection/view to match final state
[void]$ws.Range("E47").Select()
